# Include a Recovery on Non Closed Session
#
# After Excel's crash/auto-recovery, the in-progress item edit (which had
# been lost when the previous session did not close cleanly) is re-entered:
# the product row is updated from the old "SD008" item to the new "F030"
# item (name/model in column A, internal code in column B, brand in column
# P and barcode in column T all reference the same item code/name), and the
# sale price (column G) is corrected from 2 to 70.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nombre / Marca -> "F030-PANTALONETA" (was "SD008-Ajicero")
$ws.Range("A2").Value = "F030-PANTALONETA"
$ws.Range("P2").Value = "F030-PANTALONETA"

# Codigo Interno / Cod barras -> "F030" (was "SD008")
$ws.Range("B2").Value = "F030"
$ws.Range("T2").Value = "F030"

# Precio Unitario Venta -> 70 (was 2)
$ws.Range("G2").Value = 70
